$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "190 Packets" header (Z1, shared string) becomes "200 Packets"
$ws.Range("Z1").Value = "200 Packets"

# New Y/Z/AA (Acc/Loss/Time for the 200-packet trial run) for rows 3-51
$packetData = @(
    @(3, 89.105224609375, 0.25438294162742697, 73.582070589065495),
    @(4, 89.036452770233097, 0.254669853235239, 67.9497230052948),
    @(5, 89.3352627754211, 0.253604208176749, 66.092020034789996),
    @(6, 88.517087697982703, 0.26864633259369203, 66.409309387207003),
    @(7, 89.102852344512897, 0.26269290123330002, 66.010299682617102),
    @(8, 88.602459430694495, 0.26003783162621602, 66.264995574951101),
    @(9, 89.280718564987097, 0.258972979994395, 65.976788997650104),
    @(10, 88.770842552185002, 0.25462025246967701, 65.205659151077199),
    @(11, 89.340007305145207, 0.253632214114274, 65.611485958099294),
    @(12, 88.009577989578204, 0.27674202843933299, 66.539926290512),
    @(13, 89.230912923812795, 0.269551629044956, 65.167735099792395),
    @(14, 89.055424928665104, 0.25408717550116899, 65.925595998763995),
    @(15, 88.8989031314849, 0.258276793598396, 67.157522201538001),
    @(16, 88.965302705764699, 0.25737923893893799, 65.675297260284395),
    @(17, 88.922619819641099, 0.26277419243781802, 66.605241537094102),
    @(18, 88.275194168090806, 0.27176142808745202, 66.532241106033297),
    @(19, 89.043563604354802, 0.24872770649891701, 66.140627145767198),
    @(20, 89.5795285701751, 0.25570197520383497, 65.496098279953003),
    @(21, 89.304435253143296, 0.26295139850374, 66.7325856685638),
    @(22, 89.375579357147203, 0.26155262232016901, 66.541148662567096),
    @(23, 89.1194522380828, 0.257388196507287, 66.388028383255005),
    @(24, 88.417482376098604, 0.260061346537801, 65.9139981269836),
    @(25, 89.107596874237004, 0.25518282571728901, 66.501577854156494),
    @(26, 89.3210351467132, 0.25594247625935002, 65.911603450775104),
    @(27, 89.418268203735295, 0.25779701111085102, 65.8059756755828),
    @(28, 89.444351196289006, 0.24904361947571199, 66.723790645599294),
    @(29, 88.5052263736724, 0.26518077581271299, 65.204807519912706),
    @(30, 89.344745874404893, 0.25775008423846102, 66.849703073501502),
    @(31, 88.877558708190904, 0.25866306127234401, 65.678332090377793),
    @(32, 89.854627847671495, 0.24555993769185699, 66.234234571456895),
    @(33, 89.264118671417194, 0.25078757245077699, 65.366728305816594),
    @(34, 89.010363817214895, 0.26172375478293702, 65.6241614818573),
    @(35, 88.268077373504596, 0.25872530461726601, 66.216899871826101),
    @(36, 88.2467329502105, 0.26951708561742199, 66.761482715606604),
    @(37, 88.543170690536499, 0.26200680068615201, 65.874914646148596),
    @(38, 88.939219713211003, 0.25170213104571498, 66.123952150344806),
    @(39, 88.4388267993927, 0.25689226990108799, 65.218166351318303),
    @(40, 88.436454534530597, 0.26715390335801198, 66.768859863281193),
    @(41, 89.247518777847205, 0.25163316931152802, 65.841428995132404),
    @(42, 89.456212520599294, 0.251611393987623, 66.209949731826697),
    @(43, 88.386654853820801, 0.267891941957629, 65.915704965591402),
    @(44, 88.557404279708805, 0.25702410885892202, 65.625644922256399),
    @(45, 88.400882482528601, 0.259289089396744, 65.7077312469482),
    @(46, 89.093369245529104, 0.25922640642127298, 66.727711915969806),
    @(47, 88.524198532104407, 0.26364589213584899, 65.662639856338501),
    @(48, 89.226174354553194, 0.25290057101533198, 65.695250272750798),
    @(49, 88.045156002044607, 0.26422306826310699, 65.567113876342702),
    @(50, 88.984274864196706, 0.25475019325397502, 66.367227077484102),
    @(51, 88.389027118682804, 0.25856195903409301, 65.674752950668307)
)

foreach ($entry in $packetData) {
    $row = $entry[0]
    $ws.Cells.Item($row, 25).Value = $entry[1]   # column Y
    $ws.Cells.Item($row, 26).Value = $entry[2]   # column Z
    $ws.Cells.Item($row, 27).Value = $entry[3]   # column AA
}

# Match the author's final selection
$ws.Range("AA3:AA51").Select()
